# Commit: "Updated Griffin and LightInject."
#
# The underlying data table on the "Tabelle1" worksheet has a row for the
# "LightInject" DI container (row 8). Its three benchmark numbers (columns
# B/C/D) are corrected:
#   B8 (column 1 numbers): 90 -> 74
#   C8 (column 2 numbers): 84 -> 78
#   D8 (column 3 numbers): 96 -> 88
#
# Three bar3D charts on the same sheet plot those columns, so their cached
# series data for "LightInject" (category index 6) updates along with the
# cells. The author's workbook also ends up with cell B8 selected/active.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 74
$ws.Range("C8").Value = 78
$ws.Range("D8").Value = 88

# Ask Excel to recalculate everything and pull the charts' cached data back
# in line with the new cell values.
$excel.CalculateFullRebuild()
$wb.RefreshAll()

$chartObjects = $ws.ChartObjects()
for ($i = 1; $i -le $chartObjects.Count; $i++) {
    $chartObjects.Item($i).Chart.Refresh()
}

# Leave the sheet the way it ended up in the source workbook: cell B8
# selected as the active cell.
$ws.Activate() | Out-Null
$ws.Range("B8").Select() | Out-Null
